$d = $word.ActiveDocument

### Bold-only fixes (table header cells): <w:b w:val="0"/> -> <w:b/> ###

# Table 1, Row 1, Col 2: "Descrição do produto" header cell
$d.Tables.Item(1).Cell(1, 2).Range.Font.Bold = $true

# Table 2, Row 1, Col 1: "Marca" header cell
$d.Tables.Item(2).Cell(1, 1).Range.Font.Bold = $true

# Table 2, Row 1, Col 2: "Participação no mercado (%)" header cell
$d.Tables.Item(2).Cell(1, 2).Range.Font.Bold = $true

### Text replacements ###

$found = $d.Content.Find.Execute("É uma bebida versátil que pode ser apreciada quente ou fria, com ou sem leite, e com diferentes especiarias e adoçantes.", $true, $false, $false, $false, $false, $true, 1, $false, "温かくても冷たくても、ミルクの有無にかかわらず、さまざまなスパイスや甘味料と一緒に楽しめる多用途の飲み物です。", 2)
if (-not $found) { throw "Replace failed for item 0" }

$found = $d.Content.Find.Execute("Desfrute da rica e aromática experiência do Mystic Spice Premium Chai Tea, uma mistura cuidadosamente elaborada que homenageia as tradições clássicas do chai indiano.", $true, $false, $false, $false, $false, $true, 1, $false, "インドのチャイの時代を超越した伝統に敬意を表し、細心の注意を払って作られたブレンドである Mystic Spice Premium Chai Tea の豊かで香り高い抱擁をお楽しみください。", 2)
if (-not $found) { throw "Replace failed for item 1" }

# "Cada xícara oferece..." occurs twice in the document; target only the occurrence
# inside Table 1, Row 2, Col 2 (the one outside any table must stay untouched).
$cell = $d.Tables.Item(1).Cell(2, 2)
$found = $cell.Range.Find.Execute("Cada xícara oferece uma jornada encantadora pelos vibrantes cenários da Índia, trazendo uma experiência autêntica de chai diretamente para sua casa.", $false, $false, $false, $false, $false, $true, 1, $false, "各カップはインドの活気に満ちた風景を巡る魅惑的な旅を提供し、自宅で本格的なチャイ体験をお届けします。", 1)
if (-not $found) { throw "Replace failed for Cada xicara text" }

$found = $d.Content.Find.Execute("Mistura autêntica: Nosso chai é uma mistura harmoniosa de folhas de chá preto premium e uma seleção exclusiva de especiarias moídas, incluindo canela, cardamomo, cravo, gengibre e pimenta-do-reino.", $true, $false, $false, $false, $false, $true, 1, $false, "Mistura autêntica: nosso chai é uma mistura harmoniosa de folhas de chá preto de qualidade e uma seleção exclusiva de especiarias moídas, incluindo canela, cardamomo, cravo, gengibre e pimenta preta.", 2)
if (-not $found) { throw "Replace failed for item 3" }

$found = $d.Content.Find.Execute("Esta receita centenária promete um sabor autêntico e robusto em cada gole.", $true, $false, $false, $false, $false, $true, 1, $false, "この古くから伝わるレシピは、一口飲むごとに本格的でしっかりとした味わいを約束します。", 2)
if (-not $found) { throw "Replace failed for item 4" }

$found = $d.Content.Find.Execute("Ingredientes que melhoram a saúde: Cada ingrediente do chá Mystic Spice Chai é escolhido por seus benefícios naturais para a saúde.", $true, $false, $false, $false, $false, $true, 1, $false, "Ingredientes que melhoram a saúde: cada ingrediente do Mystic Spice Chai Tea é escolhido por seus benefícios naturais à saúde.", 2)
if (-not $found) { throw "Replace failed for item 5" }

$found = $d.Content.Find.Execute("Aroma e sabor ricos: O aroma quente e picante e o sabor profundo e revigorante do nosso chai fazem dele a bebida perfeita para começar o dia ou relaxar à noite.", $true, $false, $false, $false, $false, $true, 1, $false, "Aroma e sabor ricos: o aroma quente e picante e o sabor profundo e revigorante do nosso chai o tornam a bebida perfeita para começar o dia ou relaxar à noite.", 2)
if (-not $found) { throw "Replace failed for item 6" }

$found = $d.Content.Find.Execute("Os sabores são intensos, mas equilibrados, proporcionando uma experiência reconfortante e relaxante.", $true, $false, $false, $false, $false, $true, 1, $false, "風味は強烈でありながらバランスが取れており、快適で心地よい体験を生み出します。", 2)
if (-not $found) { throw "Replace failed for item 7" }

$found = $d.Content.Find.Execute("Opções versáteis de fabricação: Se você ama seu chai fumegante quente, como um chá gelado refrescante ou como um café com leite cremoso, nossa mistura é versátil o suficiente para atender a qualquer preferência.", $true, $false, $false, $false, $false, $true, 1, $false, "Opções versáteis de preparo: se você ama seu chai quente, como um chá gelado refrescante ou como um latte cremoso, nossa mistura é versátil o suficiente para atender a qualquer preferência.", 2)
if (-not $found) { throw "Replace failed for item 8" }

$found = $d.Content.Find.Execute("Instruções simples de preparo estão incluídas para ajudá-lo a saborear seu chai exatamente do jeito que você gosta.", $true, $false, $false, $false, $false, $true, 1, $false, "お好みの方法でチャイをお楽しみいただけるよう、簡単な淹れ方の説明書が付属しています。", 2)
if (-not $found) { throw "Replace failed for item 9" }

$found = $d.Content.Find.Execute("De origem sustentável: Comprometidos com a sustentabilidade, obtemos nossos ingredientes de pequenas fazendas que praticam a agricultura orgânica, garantindo não apenas a melhor qualidade, mas também o bem-estar do nosso planeta.", $true, $false, $false, $false, $false, $true, 1, $false, "Origem sustentável: comprometidos com a sustentabilidade, obtemos nossos ingredientes de fazendas de pequena escala que praticam a agricultura orgânica, garantindo não apenas a melhor qualidade, mas também o bem-estar do nosso planeta.", 2)
if (-not $found) { throw "Replace failed for item 10" }

$found = $d.Content.Find.Execute("Embalagem elegante: O Mystic Spice Chai Tea vem em embalagens ecológicas e com design lindo, tornando-o um presente ideal para os amantes do chá ou um deleite luxuoso para si mesmo.", $true, $false, $false, $false, $false, $true, 1, $false, "Embalagem elegante: o Mystic Spice Chai Tea vem em uma embalagem ecológica com um belo design, tornando-o o presente ideal para amantes de chá ou um agrado pessoal luxuoso.", 2)
if (-not $found) { throw "Replace failed for item 11" }

$found = $d.Content.Find.Execute("Garantia de Satisfação do Cliente: Nós apoiamos nosso produto e oferecemos uma garantia de satisfação.", $true, $false, $false, $false, $false, $true, 1, $false, "Garantia de satisfação do cliente: nos responsabilizamos por nosso produto e oferecemos uma garantia de satisfação.", 2)
if (-not $found) { throw "Replace failed for item 12" }

$found = $d.Content.Find.Execute("Se o Mystic Spice Chai Tea não atender suas expectativas, estamos comprometidos em resolver da melhor maneira possível.", $true, $false, $false, $false, $false, $true, 1, $false, "Mystic Spice Chai Tea がお客様のご期待に添えない場合は、当社が改善するよう努めます。", 2)
if (-not $found) { throw "Replace failed for item 13" }

$found = $d.Content.Find.Execute("Ideal para: Entusiastas do chá, indivíduos preocupados com a saúde, amantes de bebidas quentes e picantes e qualquer pessoa que queira explorar os ricos sabores do tradicional chai indiano.", $true, $false, $false, $false, $false, $true, 1, $false, "Ideal para: entusiastas do chá, indivíduos preocupados com a saúde, amantes de bebidas quentes e picantes e qualquer pessoa que queira explorar os ricos sabores do chai indiano tradicional.", 2)
if (-not $found) { throw "Replace failed for item 14" }

$found = $d.Content.Find.Execute("Tetley: A Tetley é uma empresa britânica de chá que tem forte presença na América Latina, especialmente no Brasil, onde é líder de mercado.", $true, $false, $false, $false, $false, $true, 1, $false, "Tetley: a Tetley é uma empresa britânica de chá que tem forte presença na América Latina, sobretudo no Brasil, onde é líder de mercado.", 2)
if (-not $found) { throw "Replace failed for item 15" }

$found = $d.Content.Find.Execute("Teavana: A Teavana é uma empresa de chá com sede nos EUA que pertence à Starbucks e opera em vários países da América Latina, como México, Colômbia e Peru.", $true, $false, $false, $false, $false, $true, 1, $false, "Teavana: a Teavana é uma empresa de chá com sede nos Estados Unidos, de propriedade da Starbucks e que opera em vários países da América Latina, como México, Colômbia e Peru.", 2)
if (-not $found) { throw "Replace failed for item 16" }

$found = $d.Content.Find.Execute("David's Tea: A David's Tea é uma empresa canadense de chá que está presente em alguns países da América Latina, como Chile e Costa Rica.", $true, $false, $false, $false, $false, $true, 1, $false, "David's Tea: a David's Tea é uma empresa canadense de chá que está presente em alguns países da América Latina, como Chile e Costa Rica.", 2)
if (-not $found) { throw "Replace failed for item 17" }

$found = $d.Content.Find.Execute("Marcas locais: Existem também várias marcas locais que oferecem produtos de chá Chai na América Latina, como Mate Factor, Chai Mate e Chai Brasil.", $true, $false, $false, $false, $false, $true, 1, $false, "Marcas locais: existem também várias marcas locais que oferecem produtos de chá Chai na América Latina, como Mate Factor, Chai Mate e Chai Brasil.", 2)
if (-not $found) { throw "Replace failed for item 18" }

# "Estratégia de preços" header cell (Table 2, Row 1, Col 3) also gets the bold fix
$cell = $d.Tables.Item(2).Cell(1, 3)
$cell.Range.Font.Bold = $true
$found = $cell.Range.Find.Execute("Estratégia de preços", $false, $false, $false, $false, $false, $true, 1, $false, "Estratégias de preços", 1)
if (-not $found) { throw "Replace failed for Estrategia de precos text" }

$found = $d.Content.Find.Execute("Eles desempenham um papel crucial na visibilidade e acessibilidade dos produtos de chai e podem influenciar a percepção e a preferência do consumidor.", $true, $false, $false, $false, $false, $true, 1, $false, "小売業者はチャイ ティー製品の最も目に付きやすくアクセスしやすいチャネルであり、消費者のチャイ ティー製品に対する認識、好み、購入に影響を与える可能性があります。", 2)
if (-not $found) { throw "Replace failed for item 20" }

Write-Host "All replacements completed."